# "Overskuelig kode for overblikskærmen"
# Remove the duplicate "LunaSensor / Sensor2" log rows (rows 8, 10 and 13
# were exact duplicates of row 7), add a "Btw." remark on row 4, and
# update the current selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the duplicate rows, starting from the bottom so row numbers
# above the deletion point stay valid.
$ws.Rows.Item(13).Delete()
$ws.Rows.Item(10).Delete()
$ws.Rows.Item(8).Delete()

# Add the new "Btw." remark for the borger-kræver-to-personaler entry.
$ws.Range("E4").Value = "Btw."

# Update the selected cell to match where the user ended up working.
$ws.Range("D10").Select()
